# SDDX - Xagio Work Site Building August
# Insert a new note at E4 ("I have 8 builds to do with xagio before the
# month ends"), pushing the existing three rows (xagio work / the two
# URLs) further down the column, and move the selection to the new cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remember the current (pre-edit) content of the three existing rows ---
$oldE4Text = $ws.Range("E4").Value2   # "xagio work"          (bold, blue fill)
$oldE7Text = $ws.Range("E7").Value2   # "https://treeserviceannarbor.com/wp-admin/"
$oldE8Text = $ws.Range("E8").Value2   # "lansingtreeservice.org"

# --- Move the two plain link rows down first (E7->E11, E8->E12) ---------
# These carry the default/normal style, so a straight value copy is enough.
$ws.Range("E12").Value = $oldE8Text
$ws.Range("E11").Value = $oldE7Text

# --- Move the "xagio work" row down to E8, carrying its original format -
# (bold font + themed fill). Copy/PasteSpecial(formats) preserves the
# existing theme-based fill instead of recreating it as a raw RGB fill.
$ws.Range("E4").Copy() | Out-Null
$ws.Range("E8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Application.CutCopyMode = $false
$ws.Range("E8").Value = $oldE4Text

# --- Clear the now-vacated source rows (old E4 content moved away, old
# E7 emptied since its content now lives at E11) ------------------------
$ws.Range("E7").ClearContents()

# --- Put the new note in E4 with a fresh style: bold text, no fill ------
$ws.Range("E4").ClearFormats()
$ws.Range("E4").Font.Bold = $true
$ws.Range("E4").Value = "I have 8 builds to do with xagio before the month ends"

# --- Move the active selection to the new note --------------------------
$ws.Range("E4").Select() | Out-Null
